# Update "想去人数" (interest count) figures for refreshed data pull.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 7435
$ws1.Range("F7").Value = 468
$ws1.Range("F8").Value = 184
$ws1.Range("F9").Value = 1067
$ws1.Range("F10").Value = 533
$ws1.Range("F11").Value = 16
$ws1.Range("F12").Value = 159
$ws1.Range("F13").Value = 194
$ws1.Range("F14").Value = 689

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 7435
$ws4.Range("F7").Value = 468
$ws4.Range("F8").Value = 184
$ws4.Range("F9").Value = 1067
$ws4.Range("F10").Value = 534
$ws4.Range("F11").Value = 16
$ws4.Range("F12").Value = 159
$ws4.Range("F13").Value = 194
$ws4.Range("F14").Value = 689

$wb.Save()
